# The pptx writer now uses the table's own internal column widths
# (instead of always dividing the shape width evenly across columns).
# For this deck that means the first table on slide 6 gets its two
# grid columns resized from 197pt (2,501,900 EMU) to 198pt
# (2,514,600 EMU) each -- matching the column widths already used by
# the second table on the same slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$tbl = $s.Shapes.Item(3).Table
$tbl.Columns.Item(1).Width = 198
$tbl.Columns.Item(2).Width = 198
